$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural changes -------------------------------------------------
# Insert two new (blank) columns before column F, pushing the existing
# F/G/H data out to H/I/J.
$ws.Range("F1:G1").EntireColumn.Insert()

# Insert three new (blank) rows before row 17, pushing the existing
# rows 17-19 down to 20-22.
$ws.Range("A17:A19").EntireRow.Insert()

# The row insert copies the formatting of the row above into every
# column that has data further down the sheet; remove the cells that
# must stay blank on the new rows.
$ws.Range("C17:C19").Clear()
$ws.Range("I17:J19").Clear()

# --- Move the old column C "pe4 / pe5 / pa10" values into the new
#     rows 17-19 (column B) and clear them out of column C -------------
$ws.Range("B17").Value = "pe4"
$ws.Range("B18").Value = "pe5"
$ws.Range("B19").Value = "pa10"

$ws.Range("C14:C16").Clear()

# --- New labels in column A for the new rows ---------------------------
$ws.Range("A17").Value = "Ina-2"
$ws.Range("A18").Value = "Inb-2"
$ws.Range("A19").Value = "PWM-2"

# --- New NUCLEO pinout table in columns D/E/F --------------------------
$ws.Range("D1").Value = "NUCLEO"
$ws.Range("E1").Value = "scription"

$ws.Range("D3").Value = "pf1"
$ws.Range("E3").Value = "CN9-19"

$ws.Range("D4").Value = "pf0"
$ws.Range("E4").Value = "CN9-21"

$ws.Range("D11").Value = "i2c_2"

$ws.Range("D13").Value = "1"
$ws.Range("E13").Value = "2"

$ws.Range("D14").Value = "PE4"
$ws.Range("E14").Value = "CN9-16"

$ws.Range("D15").Value = "PE5"
$ws.Range("E15").Value = "CN9-18"

$ws.Range("D16").Value = "PB10"
$ws.Range("E16").Value = "CN10-32/"
$ws.Range("F16").Value = "TIM2_CH3"

$ws.Range("D17").Value = "PE6"
$ws.Range("E17").Value = "CN9-20"

$ws.Range("D18").Value = "PE3"
$ws.Range("E18").Value = "CN9-22"

$ws.Range("D19").Value = "PB11"
$ws.Range("E19").Value = "CN10-34/"
$ws.Range("F19").Value = "TIM2_CH4"

# --- Apply the same plain-text cell format used throughout the sheet
#     (numFmtId 49 / style index 1) to every freshly populated cell ----
$newCells = @("D1","E1","D3","E3","D4","E4","D11","D13","E13","D14","E14","D15","E15","D16","E16","F16","D17","E17","D18","E18","D19","E19","F19")
foreach ($addr in $newCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Restore the selection shown in the saved workbook -----------------
$ws.Range("F20").Select()
